$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-18.82065889985288,2.594480925474296,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,2.479633885105085,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,2.436607502240222,-18.82065889985288,-18.82065889985288,2.48341663565457,-18.82065889985288,1.871974478781928,-18.82065889985288,2.393397514088249,-18.82065889985288),
    @(-18.82065889985288,0.4441066761610361,-18.82065889985288,-18.82065889985288,-18.82065889985288,2.110116690655433,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(2.820968214570909,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,2.773070041465228,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(3.693051590900136,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,1.524634801558814,-18.82065889985288,2.259167647087663),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,2.280437714871896,-18.82065889985288,2.566180857864985,-18.82065889985288,-18.82065889985288,-18.82065889985288,1.33815908569216),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,1.508494704960074,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,2.324014900336324,1.614971688052147),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,2.082205633409901),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,2.284293630348477,-18.82065889985288),
    @(-18.82065889985288,1.466468393316565,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,0.4443814922296981,0.8759707946693541,1.186402596029484,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,-18.82065889985288,0.4252488844938448,1.033474425829998,1.373897186954315,-18.82065889985288),
    @(-18.82065889985288,-18.82065889985288,4.321925134723702,-18.82065889985288,-18.82065889985288,-18.82065889985288,1.826294383889045,2.04091345848125,-18.82065889985288,-18.82065889985288),
    @(-18.82065889985288,1.210354354185747,-18.82065889985288,-18.82065889985288,3.848773999057811,-18.82065889985288,2.214656381267064,1.82701470759852,-18.82065889985288,2.429464264712927),
    @(-18.82065889985288,1.076366690699287,-18.82065889985288,2.449157783708481,-18.82065889985288,3.286926315709758,2.445764510733449,-18.82065889985288,-18.82065889985288,-18.82065889985288)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $row[$j]
    }
}